# Apply the edits described by the commit "completed summary dataframe and timings"
# 1) Update the narrative shared strings in column F (rows 2-14) with the new
#    assessment descriptions.
# 2) Update the "Assessment Duration" values in column C for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update the assessment description text (column F, rows 2-14) ---

$ws.Range("F2").Value = "Practical exercises provide a direct demonstration of the ability to analyze release components. Assessing the practical application of identifying and understanding the different components of a release.`nA1: Analyse release components (ICT-DIT-3014-1.1)"

$ws.Range("F3").Value = "Practical exercises simulate real-world scenarios requiring coordination. Assessing the ability to communicate effectively and align release activities with stakeholders' needs and timelines.`nA2: Coordinate with relevant stakeholders on release scheduling to align release processes and procedures (ICT-DIT-3014-1.1)"

$ws.Range("F4").Value = "Practical exercises assess the application of knowledge regarding scripts and tools, and their usage in software integration and deployment. It tests understanding of when and how to use specific tools.`nK1: Types and usage of scripts and tools for integrating and deploying software products (ICT-DIT-3014-1.1)"

$ws.Range("F5").Value = "Practical assessment demonstrates the ability to choose appropriate scripts and tools. Focusing on effectiveness and suitability for specific integration and deployment tasks, evaluating selection skills.`nA3: Select appropriate scripts and tools for integrating and deploying software products (ICT-DIT-3014-1.1)"

$ws.Range("F6").Value = "Practical application tests the understanding of software configuration procedures by requiring candidates to execute these procedures, showing they understand each step's purpose and sequence.`nK2: Software configuration procedures (ICT-DIT-3014-1.1)"

$ws.Range("F7").Value = "Configuration tests are included to check the candidate knows the configuration tests performed, and their purpose in ensuring the reliability and stability of the software product.`nK3: Configuration tests and their purposes (ICT-DIT-3014-1.1)"

$ws.Range("F8").Value = "Practical exercises test the ability to configure software. It assesses the ability to translate theoretical knowledge into practical skills in diverse platform configurations.`nA4: Configure software products to integrate and deploy software releases to various platforms (ICT-DIT-3014-1.1)"

$ws.Range("F9").Value = "Practical activities demonstrate competency in executing configuration tests. Assessing the ability to follow procedures accurately and efficiently, while understanding the expected outcomes.`nA5: Execute configuration tests on platform specific versions of software products in line with testing procedures (ICT-DIT-3014-1.1)"

$ws.Range("F10").Value = "Practical activities ensures the candidate can determine the outcome, issues or successes of testing procedures, and ensures the software is functional after configuration tests.`nK4: Interpretation of configuration test results (ICT-DIT-3014-1.1)"

$ws.Range("F11").Value = "Practical tests involve identifying the cause and remedy to issues, which allows for assessment of their ability to diagnose effectively. The configuration allows a setting for practical investigation.`nA6: Diagnose issues surfaced from configuration testing (ICT-DIT-3014-1.1)"

$ws.Range("F12").Value = "The practical element ensures the candidate understands the software configuration and deployment process, and they have appropriate knowledge to perform actions that will be tested in action.`nK5: Elements of the software configuration and deployment process (ICT-DIT-3014-1.1)"

$ws.Range("F13").Value = "Practical implementation allows for the ability to determine improvements from the outcome of testing, and their experience, which allows for better modification of software code.`nA7: Identify potential improvements and modifications to the software configuration and deployment process or the software code (ICT-DIT-3014-1.1)"

$ws.Range("F14").Value = "The practical element ensures the candidate can properly modify platform-specific software products and processes, while retaining previous functionality and applying improvements.`nA8: Implement modifications to platform-specific software products and processes (ICT-DIT-3014-1.1)"

# Restore automatic row heights so the newline characters inserted above
# don't leave a stray explicit/custom row height behind.
$ws.Rows("2:14").AutoFit()

# --- 2) Update Assessment Duration values (column C) ---

$ws.Range("C2").Value = 5
$ws.Range("C3").Value = 5
$ws.Range("C4").Value = 10
$ws.Range("C6").Value = 10
$ws.Range("C7").Value = 10
$ws.Range("C10").Value = 10
$ws.Range("C12").Value = 10
